# BD.xlsx edit: add "nombre" and "apellido" columns to the "usuario" sheet
# documentation table, and leave the "usuario" sheet active (matching the
# tab/selection state recorded in the workbook).

$wb = $excel.ActiveWorkbook

# The previous selection state had the "estado_aut" sheet active with
# cell B3 selected; after this edit it ends up with E9 selected (but no
# longer the active tab).
$wsEstadoAut = $wb.Worksheets.Item("estado_aut")
$wsEstadoAut.Select()
$wsEstadoAut.Range("E9").Select()

# Work on the "usuario" sheet: insert two new documentation rows (for the
# "nombre" and "apellido" columns of the usuario table) right after the
# existing "id_usuario" row, pushing the "id_rol" / "id_estado_aut" rows
# down from rows 3-4 to rows 5-6.
$ws = $wb.Worksheets.Item("usuario")
$ws.Select()

$ws.Rows("3:4").Insert()

$ws.Range("A3").Value = "nombre"
$ws.Range("B3").Value = "VARCHAR(50)"
$ws.Range("C3").Value = "NOT NULL"
$ws.Range("E3").Value = "nombre del usuario"

$ws.Range("A4").Value = "apellido"
$ws.Range("B4").Value = "VARCHAR(50)"
$ws.Range("C4").Value = "NOT NULL"
$ws.Range("E4").Value = "apellido del usuario"

# Final selection on the "usuario" sheet (now the active tab).
$ws.Range("E5").Select()
